$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cell values taken from the source diff (coin price / volume refresh,
# plus two pairs of rows - Dai/ShibaInu and FraxShare/Aave - that swapped
# rank order and so swapped their B/C/D/E contents).
$updates = @(
    @{ Ref = 'D2'; Value = '30.379.22' },
    @{ Ref = 'E2'; Value = '  -0.77%  ' },
    @{ Ref = 'D3'; Value = '1.871.42' },
    @{ Ref = 'E3'; Value = '  -0.35%  ' },
    @{ Ref = 'D4'; Value = '1.0000' },
    @{ Ref = 'E4'; Value = '  -0.05%  ' },
    @{ Ref = 'D5'; Value = '243.98' },
    @{ Ref = 'E5'; Value = '  -1.70%  ' },
    @{ Ref = 'E6'; Value = '  -0.01%  ' },
    @{ Ref = 'D7'; Value = '0.4712' },
    @{ Ref = 'E7'; Value = '  -1.03%  ' },
    @{ Ref = 'D8'; Value = '0.2886' },
    @{ Ref = 'E8'; Value = '  -1.35%  ' },
    @{ Ref = 'D9'; Value = '0.06470' },
    @{ Ref = 'E9'; Value = '  -0.87%  ' },
    @{ Ref = 'D10'; Value = '21.98' },
    @{ Ref = 'E10'; Value = '  -0.09%  ' },
    @{ Ref = 'D11'; Value = '0.07793' },
    @{ Ref = 'E11'; Value = '  +0.72%  ' },
    @{ Ref = 'D12'; Value = '96.22' },
    @{ Ref = 'E12'; Value = '  -0.61%  ' },
    @{ Ref = 'D13'; Value = '1.867.10' },
    @{ Ref = 'E13'; Value = '  -0.57%  ' },
    @{ Ref = 'D14'; Value = '0.7246' },
    @{ Ref = 'E14'; Value = '  -2.13%  ' },
    @{ Ref = 'D15'; Value = '5.142' },
    @{ Ref = 'E15'; Value = '  -1.26%  ' },
    @{ Ref = 'D16'; Value = '282.31' },
    @{ Ref = 'E16'; Value = '  +2.94%  ' },
    @{ Ref = 'D17'; Value = '30.369.18' },
    @{ Ref = 'E17'; Value = '  -1.14%  ' },
    @{ Ref = 'E18'; Value = '  -1.73%  ' },
    @{ Ref = 'B19'; Value = 'ShibaInu' },
    @{ Ref = 'C19'; Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib' },
    @{ Ref = 'D19'; Value = '0.000007511' },
    @{ Ref = 'E19'; Value = '  -0.27%  ' },
    @{ Ref = 'B20'; Value = 'Dai' },
    @{ Ref = 'C20'; Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai' },
    @{ Ref = 'D20'; Value = '1.0000' },
    @{ Ref = 'E20'; Value = '  -0.04%  ' },
    @{ Ref = 'D21'; Value = '2.113.93' },
    @{ Ref = 'E21'; Value = '  -0.43%  ' },
    @{ Ref = 'D22'; Value = '0.9999' },
    @{ Ref = 'E22'; Value = '  -0.09%  ' },
    @{ Ref = 'D23'; Value = '5.263' },
    @{ Ref = 'E23'; Value = '  +0.13%  ' },
    @{ Ref = 'D24'; Value = '6.245' },
    @{ Ref = 'E24'; Value = '  +0.58%  ' },
    @{ Ref = 'D25'; Value = '163.66' },
    @{ Ref = 'E25'; Value = '  -0.99%  ' },
    @{ Ref = 'D26'; Value = '9.063' },
    @{ Ref = 'E26'; Value = '  -1.49%  ' },
    @{ Ref = 'D27'; Value = '18.75' },
    @{ Ref = 'E27'; Value = '  -0.68%  ' },
    @{ Ref = 'D28'; Value = '1.883' },
    @{ Ref = 'E28'; Value = '  -1.76%  ' },
    @{ Ref = 'E29'; Value = '  -1.39%  ' },
    @{ Ref = 'D30'; Value = '0.09623' },
    @{ Ref = 'E30'; Value = '  -2.31%  ' },
    @{ Ref = 'D31'; Value = '1.486' },
    @{ Ref = 'E31'; Value = '  -1.09%  ' },
    @{ Ref = 'D32'; Value = '4.235' },
    @{ Ref = 'E32'; Value = '  -1.29%  ' },
    @{ Ref = 'D33'; Value = '4.118' },
    @{ Ref = 'E33'; Value = '  +0.02%  ' },
    @{ Ref = 'D34'; Value = '0.04821' },
    @{ Ref = 'E34'; Value = '  -0.30%  ' },
    @{ Ref = 'D35'; Value = '1.122' },
    @{ Ref = 'E35'; Value = '  -0.39%  ' },
    @{ Ref = 'D36'; Value = '0.6904' },
    @{ Ref = 'E36'; Value = '  -0.82%  ' },
    @{ Ref = 'D37'; Value = '2.715' },
    @{ Ref = 'E37'; Value = '  -0.07%  ' },
    @{ Ref = 'D38'; Value = '0.01892' },
    @{ Ref = 'E38'; Value = '  +0.69%  ' },
    @{ Ref = 'D39'; Value = '2.821' },
    @{ Ref = 'E39'; Value = '  +2.07%  ' },
    @{ Ref = 'B40'; Value = 'Aave' },
    @{ Ref = 'C40'; Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave' },
    @{ Ref = 'D40'; Value = '75.46' },
    @{ Ref = 'E40'; Value = '  +2.75%  ' },
    @{ Ref = 'B41'; Value = 'FraxShare' },
    @{ Ref = 'C41'; Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs' },
    @{ Ref = 'D41'; Value = '6.246' },
    @{ Ref = 'E41'; Value = '  -0.41%  ' },
    @{ Ref = 'D42'; Value = '0.4232' },
    @{ Ref = 'E42'; Value = '  -0.26%  ' },
    @{ Ref = 'D43'; Value = '1.933' },
    @{ Ref = 'E43'; Value = '  -3.11%  ' },
    @{ Ref = 'E44'; Value = '  -0.10%  ' },
    @{ Ref = 'D45'; Value = '0.8280' },
    @{ Ref = 'E45'; Value = '  -1.17%  ' },
    @{ Ref = 'D46'; Value = '100.94' },
    @{ Ref = 'E46'; Value = '  -1.22%  ' },
    @{ Ref = 'D47'; Value = '9.675' },
    @{ Ref = 'E47'; Value = '  +3.26%  ' },
    @{ Ref = 'D48'; Value = '35.32' },
    @{ Ref = 'E48'; Value = '  -0.40%  ' },
    @{ Ref = 'D49'; Value = '6.974' },
    @{ Ref = 'E49'; Value = '  -0.90%  ' },
    @{ Ref = 'D50'; Value = '902.04' },
    @{ Ref = 'E50'; Value = '  -0.96%  ' },
    @{ Ref = 'D51'; Value = '0.05728' },
    @{ Ref = 'E51'; Value = '  +0.72%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Ref)
    # Many "Price" column values are plain digits with dots (e.g. "1.0000",
    # "6.974") which Excel's COM layer would otherwise auto-coerce into a
    # Number (dropping trailing zeros / switching representation). Detect
    # that case, force Text format for the write, then restore the default
    # "Normal" style so the cell's style index is left exactly as it was.
    if ($u.Value -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
